# Automatische test-sync: 2025-08-05 18:30:50
# Appends a new test-mail log row (row 33) to the "Logs" sheet and
# refreshes the dependent "Dashboard" summary + conditional formatting
# ranges so they keep covering the full data range.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Row 33: new test-mail entry (columns A..J -> 1..10)
$logs.Cells.Item(33, 1).Value = "Ik heb nog geen geld terug."
$logs.Cells.Item(33, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(33, 3).Value = "Testmail #12: Ik heb nog geen geld terug."
$logs.Cells.Item(33, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(33, 5).Value = "Beste klant,
Bedankt voor uw bericht. Ik begrijp dat u nog geen geld hebt ontvangen en ik help u hier graag mee verder. Om uw vraag goed te kunnen beantwoorden, zou ik wat meer informatie nodig hebben. Kunt u mij laten weten om welke transactie het gaat en eventueel het bijbehorende referentienummer? Op die manier kan ik het voor u nakijken en u verder helpen.
Ik kijk uit naar uw reactie.
Met vriendelijke groet,
[Naam]
E-mailassistent"
$logs.Cells.Item(33, 6).Value = "2025-08-05 18:30:33"
$logs.Cells.Item(33, 7).Value = "Ja"
$logs.Cells.Item(33, 8).Value = "Nee"
$logs.Cells.Item(33, 9).Value = "Ja"
$logs.Cells.Item(33, 10).Value = "Ja"

# Extend the conditional-formatting ranges (D/G/H/I/J) from row 32 to the
# newly-added row 33, matching the widened dimension A1:J33.
$logs.Range("D2:D32").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D33"))
$logs.Range("G2:G32").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G33"))
$logs.Range("H2:H32").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H33"))
$logs.Range("I2:I32").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I33"))
$logs.Range("J2:J32").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J33"))

# Refresh the Dashboard roll-up: one more "Retour / Terugbetaling" entry.
$dash.Range("B5").Value = 3
